$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $ws.Range("Z1").Copy()
    $cell.PasteSpecial(-4122)
}

$ws.Range("D2").Value = '39.931.26'
$ws.Range("E2").Value = '  +0.93%  '
$ws.Range("D3").Value = '2.224.09'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  -0.02%  '
Set-TextValue $ws.Range("D5") '292.90'
$ws.Range("E5").Value = '  -1.59%  '
Set-TextValue $ws.Range("D6") '87.82'
$ws.Range("E6").Value = '  +6.05%  '
Set-TextValue $ws.Range("D7") '0.517'
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("E8").Value = '  -0.05%  '
Set-TextValue $ws.Range("D9") '0.474'
$ws.Range("E9").Value = '  +0.23%  '
Set-TextValue $ws.Range("D10") '30.66'
$ws.Range("E10").Value = '  +1.96%  '
Set-TextValue $ws.Range("D11") '0.0787'
$ws.Range("E11").Value = '  +1.06%  '
Set-TextValue $ws.Range("D12") '47.59'
$ws.Range("E12").Value = '  +1.76%  '
Set-TextValue $ws.Range("D13") '0.109'
$ws.Range("E13").Value = '  +1.55%  '
Set-TextValue $ws.Range("D14") '6.44'
$ws.Range("E14").Value = '  +2.08%  '
$ws.Range("D15").Value = '2.565.42'
$ws.Range("E15").Value = '  -0.21%  '
Set-TextValue $ws.Range("D16") '14.13'
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").Value = '2.223.77'
Set-TextValue $ws.Range("D18") '0.733'
$ws.Range("E18").Value = '  +1.87%  '
$ws.Range("D19").Value = '39.876.06'
$ws.Range("E19").Value = '  +0.97%  '
Set-TextValue $ws.Range("D20") '11.53'
$ws.Range("E20").Value = '  +11.39%  '
$ws.Range("D21").Value = '0.0₃0887'
$ws.Range("E21").Value = '  +1.05%  '
Set-TextValue $ws.Range("D22") '5.86'
$ws.Range("E22").Value = '  +1.23%  '
Set-TextValue $ws.Range("D23") '65.93'
$ws.Range("E23").Value = '  +1.18%  '
Set-TextValue $ws.Range("D24") '236.20'
$ws.Range("E24").Value = '  +3.13%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  +2.17%  '
$ws.Range("E27").Value = '  +0.56%  '
Set-TextValue $ws.Range("D28") '22.86'
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("E29").Value = '  +1.11%  '
Set-TextValue $ws.Range("D30") '9.29'
$ws.Range("E30").Value = '  +1.17%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D31") '32.89'
$ws.Range("E31").Value = '  +1.60%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D32") '152.85'
$ws.Range("E32").Value = '  +2.66%  '
$ws.Range("E33").Value = '  -0.11%  '
Set-TextValue $ws.Range("D34") '4.96'
$ws.Range("E34").Value = '  +2.20%  '
Set-TextValue $ws.Range("D35") '0.0721'
$ws.Range("E35").Value = '  +2.80%  '
$ws.Range("E36").Value = '  +1.01%  '
Set-TextValue $ws.Range("D37") '2.83'
$ws.Range("E37").Value = '  +6.29%  '
$ws.Range("E38").Value = '  +1.14%  '
Set-TextValue $ws.Range("D39") '16.01'
$ws.Range("E39").Value = '  +0.50%  '
Set-TextValue $ws.Range("D40") '0.0996'
$ws.Range("E40").Value = '  +2.66%  '
Set-TextValue $ws.Range("D41") '1.72'
$ws.Range("E41").Value = '  +2.45%  '
$ws.Range("D42").Value = '2.103.73'
$ws.Range("E42").Value = '  +9.50%  '
Set-TextValue $ws.Range("D43") '3.81'
$ws.Range("E43").Value = '  +3.12%  '
Set-TextValue $ws.Range("D44") '2.18'
$ws.Range("E44").Value = '  +6.24%  '
Set-TextValue $ws.Range("D45") '0.0271'
$ws.Range("E45").Value = '  +2.58%  '
Set-TextValue $ws.Range("D46") '10.04'
$ws.Range("E46").Value = '  +9.68%  '
Set-TextValue $ws.Range("D47") '17.69'
$ws.Range("E47").Value = '  +7.08%  '
Set-TextValue $ws.Range("D48") '2.67'
$ws.Range("E48").Value = '  +2.22%  '
$ws.Range("D49").Value = '2.436.98'
$ws.Range("E49").Value = '  -0.19%  '
Set-TextValue $ws.Range("D50") '71.12'
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D51") '1.46'
$ws.Range("E51").Value = '  +6.47%  '
